$d = $word.ActiveDocument
$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# ---------------------------------------------------------------------------
# Move the "Background" heading (blank bold paragraph + "Background"
# paragraph) so it immediately follows the Introduction section, and add the
# new Background body paragraph (tab + Iocaine Powder algorithm description).
# The trailing "_GoBack" bookmark moves along with it, to the end of the
# newly-typed body paragraph. This single range-replace both relocates the
# heading and removes it from its old spot (the source range spans the
# Introduction's last paragraph through the old "Background" paragraph).
# ---------------------------------------------------------------------------

$introPara = $d.Paragraphs(9)
$backgroundHeadingPara = $d.Paragraphs(11)

$moveRange = $d.Range($introPara.Range.Start, $backgroundHeadingPara.Range.End)

$newSectionXml = @"
<w:p $wns><w:pPr><w:ind w:firstLine="720"/></w:pPr><w:r><w:t>Our strategy performed significantly better than the Nash Equilibrium when pitted against other programs from the class. Our program finished 2</w:t></w:r><w:r><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t>nd</w:t></w:r><w:r><w:t xml:space="preserve"> overall in the preliminary tournament, scoring a total of 36 out of 44 possible points, where 2 points are awarded for a win, 1 for a tie, and 0 for a loss.</w:t></w:r></w:p><w:p $wns><w:pPr><w:rPr><w:b/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr></w:p><w:p $wns><w:pPr><w:rPr><w:b/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>Background</w:t></w:r></w:p><w:p $wns><w:r><w:rPr><w:b/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:tab/></w:r><w:r><w:t>Our rock, paper, scissors bot implemented a famous algorithm, Iocaine Powder. First developed by Dan Egnor to compete in the First International RoShamBo Programming Competition, the bot employs a mixture of strategies to both exploit the opponent and remain unexploitable. The algorithm chooses between playing a random move, a move based on an analysis of the frequency of the opponent&#8217;s moves, or a move based on an analysis of the history of the opponent&#8217;s moves.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
"@

$moveRange.InsertXML($newSectionXml)

# ---------------------------------------------------------------------------
# Add a lastRenderedPageBreak before "Conclusion". After the insertion above,
# the document now has 18 paragraphs and "Conclusion" is the last one.
# ---------------------------------------------------------------------------

$conclusionPara = $d.Paragraphs($d.Paragraphs.Count)

$conclusionXml = "<w:p $wns><w:pPr><w:rPr><w:b/><w:sz w:val=`"32`"/><w:szCs w:val=`"32`"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val=`"32`"/><w:szCs w:val=`"32`"/></w:rPr><w:lastRenderedPageBreak/><w:t>Conclusion</w:t></w:r></w:p>"
$conclusionPara.Range.InsertXML($conclusionXml)
